$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in week 11 (row 12) contribution values
$ws.Range("C12").Value = 0.22700000000000001
$ws.Range("D12").Value = 0.10666666666666599
$ws.Range("E12").Value = 0.22666666666666599
$ws.Range("F12").Value = 0.16666666666666599
$ws.Range("G12").Value = 0.10666666666666599
$ws.Range("H12").Value = 0.16666666666666599

# Fill in week 12 (row 13) contribution values
$ws.Range("C13").Value = 0.16666666666666599
$ws.Range("D13").Value = 0.16666666666666599
$ws.Range("E13").Value = 0.16666666666666599
$ws.Range("F13").Value = 0.16666666666666599
$ws.Range("G13").Value = 0.16666666666666599
$ws.Range("H13").Value = 0.16666666666666599
